$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row (row => column => value)
$data = @{
    2  = @{ E=3; G=142.9073533333333; H=428.72206;         I=0.5576664151504187; J=0.5576664151504188; K=3; M=9.827738000000002;  N=29.483214;  O=0.3869625527756497; P=0.3869625527756497; Q=1404.456026833427; R=12640.10424150084; S=0.2157960196038513;  T=0.2157960196038513  }
    3  = @{ E=3; G=142.9073533333333; H=428.72206;         I=0.5576664151504187; J=0.5576664151504188; K=3; M=9.362736333333332;  N=28.088209;  O=0.368653331266258;  P=0.368653331266258;  Q=1338.003869354504; R=12042.03482419054; S=0.2055855816805139;  T=0.2055855816805139  }
    4  = @{ E=3; G=142.9073533333333; H=428.72206;         I=0.5576664151504187; J=0.5576664151504188; K=3; M=6.206655;            N=18.619965;  O=0.2443841159580923; P=0.2443841159580923; Q=886.9766391031001;  R=7982.7897519279;   S=0.1362848138660535;  T=0.1362848138660536  }
    5  = @{ E=3; G=63.967809;          H=191.903427;        I=0.2496211559306514; J=0.2496211559306514; K=3; M=9.827738000000002;  N=29.483214;  O=0.3869625527756497; P=0.3869625527756497; Q=628.6588672860421; R=5657.929805574378; S=0.09659403972573335; T=0.09659403972573337 }
    6  = @{ E=3; G=63.967809;          H=191.903427;        I=0.2496211559306514; J=0.2496211559306514; K=3; M=9.362736333333332;  N=28.088209;  O=0.368653331266258;  P=0.368653331266258;  Q=598.9137294880269; R=5390.223565392243; S=0.09202367068836866; T=0.09202367068836868 }
    7  = @{ E=3; G=63.967809;          H=191.903427;        I=0.2496211559306514; J=0.2496211559306514; K=3; M=6.206655;            N=18.619965;  O=0.2443841159580923; P=0.2443841159580923; Q=397.026121568895;  R=3573.235094120055; S=0.06100344551654933; T=0.06100344551654934 }
    8  = @{ E=3; G=49.38440333333333;  H=148.15321;         I=0.1927124289189298; J=0.1927124289189298; K=3; M=9.827738000000002;  N=29.483214;  O=0.3869625527756497; P=0.3869625527756497; Q=485.3369772463267; R=4368.032795216941; S=0.07457249344606501; T=0.07457249344606502 }
    9  = @{ E=3; G=49.38440333333333;  H=148.15321;         I=0.1927124289189298; J=0.1927124289189298; K=3; M=9.362736333333332;  N=28.088209;  O=0.368653331266258;  P=0.368653331266258;  Q=462.3731473889877; R=4161.35832650089;  S=0.07104407889737543; T=0.07104407889737545 }
    10 = @{ E=3; G=49.38440333333333;  H=148.15321;         I=0.1927124289189298; J=0.1927124289189298; K=3; M=6.206655;            N=18.619965;  O=0.2443841159580923; P=0.2443841159580923; Q=306.51195387085;   R=2758.60758483765;  S=0.04709585657548936; T=0.04709585657548936 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
